# Re-run of the backward-elimination export: the underlying Python script
# now wraps the "write results to Excel" step in a try/except, and was
# re-executed, so every per-step OLS summary block embedded in column B
# got a fresh "Date:" / "Time:" stamp from statsmodels. Propagate the new
# timestamp into every worksheet's summary text (cell B2), leaving the
# rest of each block (which is identical across the run) untouched.

$wb = $excel.ActiveWorkbook

$oldDate = "Sun, 29 Dec 2019"
$newDate = "Wed, 01 Jan 2020"
$oldTimes = @("16:11:18", "16:11:19")
$newTime = "23:18:56"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value()
    if ($null -eq $text) { continue }

    $updated = $text

    if ($updated.Contains($oldDate)) {
        $updated = $updated.Replace($oldDate, $newDate)
    }

    foreach ($oldTime in $oldTimes) {
        if ($updated.Contains($oldTime)) {
            $updated = $updated.Replace($oldTime, $newTime)
        }
    }

    if ($updated -ne $text) {
        $cell.Value = $updated
    }
}
